# Updated BGR model - 2025-08-21 09:18
$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": swap the two timeslice-order strings ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "WaD,RaD,FaD,WaP,SaD,RaP,FaP,SaP"
$wsEv.Range("C14").Value = "RaP,RaN,FaP,SaP,FaN,SaN,WaN,WaP"

# --- Sheet "re_profiles": re-order the hydro seasonal availability rows (M4:O7) ---
$wsRe = $wb.Worksheets.Item("re_profiles")
$wsRe.Range("M4").Value = "R"
$wsRe.Range("N4").Value = 0.4005462988254575
$wsRe.Range("M5").Value = "F"
$wsRe.Range("N5").Value = 0.1858508604206501
$wsRe.Range("M6").Value = "S"
$wsRe.Range("N6").Value = 0.34121824638077031
$wsRe.Range("M7").Value = "W"
$wsRe.Range("N7").Value = 0.27238459437312212
